$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.841.98"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.377.20"
$ws.Range("E3").Value = "  -3.80%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.13"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.07"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("E8").Value = "  -8.50%  "
$ws.Range("D9").Value = "2.376.33"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.34"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.47"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "2.805.29"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "60.789.82"
$ws.Range("D18").Value = "2.373.46"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.55"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.98"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "2.492.64"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "0.0₃0924"
$ws.Range("E29").Value = "  -6.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "517.60"
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("E31").Value = "  -4.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.47"
$ws.Range("E37").Value = "  -7.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.64"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.08"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.45"
$ws.Range("E43").Value = "  -6.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.24"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.34"
$ws.Range("E46").Value = "  -5.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.34"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0912"
$ws.Range("E51").Value = "  -3.08%  "
